$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: update the "Conversión del día" note text ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$text = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.07 = 7465.85 pesos`n✅ 7465.85 pesos = 2.06 = 960.81 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $text

# --- tasas!N10, O10, N12, O12: update rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 482.999
$ws2.Range("O10").Value = 3606
$ws2.Range("N12").Value = 3621
$ws2.Range("O12").Value = 466
